# Localization status report refresh ("Generate Report for Archive"):
#   - The handed-off items have moved on from "Ready for handoff" to
#     "In Translation", so every cell carrying that status text is updated.
#   - The "Status" column (and the twin status columns on the Overview
#     sheet) is re-autofit to the now-shorter text, narrowing it.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # NOTE: keep the string literal on the left of -eq. Some cells hold
        # native booleans (e.g. the "True"/"False" status columns) and
        # PowerShell's -eq coerces the right-hand side to the left operand's
        # type, so "$boolCell.Value2 -eq $oldStatus" would coerce the
        # non-empty $oldStatus string to $true and false-match every
        # truthy boolean cell.
        if ($oldStatus -eq $cell.Value2) {
            $cell.Value = $newStatus
            # Narrow the column to fit the shorter replacement text, just
            # like Excel does when you edit a cell and re-fit the column.
            $cell.EntireColumn.ColumnWidth = 12.5
        }
    }
}
